$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row.
# The refresh bumps this date by one day (45178 -> 45179) for every row
# from row 2 through row 262.
$lastRow = 262
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45179
}
